# Fruta / hortaliza, semanal
# Insert a new weekly record as row 25 in the data table, shifting the
# existing rows 25-57 down to rows 26-58.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 25. This pushes the old
# row 25 (and everything below it) down by one row.
$ws.Rows.Item(25).Insert()

# Populate the newly inserted row 25 with the new weekly price record.
$ws.Cells.Item(25, 1).Value2  = 11
$ws.Cells.Item(25, 2).Value2  = "Vega Monumental Concepción"
$ws.Cells.Item(25, 3).Value2  = "Bíobío"
$ws.Cells.Item(25, 4).Value2  = 45219
$ws.Cells.Item(25, 5).Value2  = 8
$ws.Cells.Item(25, 6).Value2  = "Fruta"
$ws.Cells.Item(25, 7).Value2  = 100107
$ws.Cells.Item(25, 8).Value2  = "Otros"
$ws.Cells.Item(25, 9).Value2  = 100107002
$ws.Cells.Item(25, 10).Value2 = "Chirimoya"
$ws.Cells.Item(25, 11).Value2 = "Cultivar IV Región"
$ws.Cells.Item(25, 12).Value2 = "Primera"
$ws.Cells.Item(25, 13).Value2 = 140
$ws.Cells.Item(25, 14).Value2 = 20000
$ws.Cells.Item(25, 15).Value2 = 21000
$ws.Cells.Item(25, 16).Value2 = 20429
$ws.Cells.Item(25, 17).Value2 = "$/bandeja 10 kilos"
$ws.Cells.Item(25, 18).Value2 = "Provincia de Limarí"
$ws.Cells.Item(25, 19).Value2 = 2043
$ws.Cells.Item(25, 20).Value2 = 10

# Keep the date column formatted the same way as the rest of column D.
$ws.Cells.Item(25, 4).NumberFormat = $ws.Cells.Item(26, 4).NumberFormat
